$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.942.21'
$ws.Range("E2").Value = '  -1.16%  '

# r3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.432.95'
$ws.Range("E3").Value = '  -1.91%  '

# r4
$ws.Range("E4").Value = '  -0.03%  '

# r5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.78'
$ws.Range("E5").Value = '  -1.24%  '

# r6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.16'
$ws.Range("E6").Value = '  -1.51%  '

# r7
$ws.Range("E7").Value = '  +0.01%  '

# r8
$ws.Range("E8").Value = '  -3.04%  '

# r9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.434.62'
$ws.Range("E9").Value = '  -2.03%  '

# r10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.23'
$ws.Range("E10").Value = '  -1.22%  '

# r12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.448'
$ws.Range("E12").Value = '  +0.50%  '

# r13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.025.11'
$ws.Range("E13").Value = '  -2.07%  '

# r14
$ws.Range("E14").Value = '  -0.28%  '

# r15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000188'
$ws.Range("E15").Value = '  -3.54%  '

# r16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.80'
$ws.Range("E16").Value = '  -3.31%  '

# r17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.948.97'
$ws.Range("E17").Value = '  -1.18%  '

# r18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.407.22'
$ws.Range("E18").Value = '  -2.15%  '

# r19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.38'
$ws.Range("E19").Value = '  -1.37%  '

# r20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.86'
$ws.Range("E20").Value = '  -3.27%  '

# r21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '381.73'
$ws.Range("E21").Value = '  -2.39%  '

# r22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.97'
$ws.Range("E22").Value = '  -3.80%  '

# r23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.548'
$ws.Range("E23").Value = '  -0.83%  '

# r24
$ws.Range("E24").Value = '  +0.00%  '

# r25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.13'
$ws.Range("E25").Value = '  -1.77%  '

# r26
$ws.Range("E26").Value = '  -5.09%  '

# r27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.93'
$ws.Range("E27").Value = '  -2.25%  '

# r28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.178'
$ws.Range("E28").Value = '  -0.76%  '

# r29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.04%  '

# r30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.47'
$ws.Range("E30").Value = '  +0.91%  '

# r31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.16'
$ws.Range("E31").Value = '  -2.80%  '

# r32
$ws.Range("E32").Value = '  -2.64%  '

# r33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.30'
$ws.Range("E33").Value = '  -1.79%  '

# r34
$ws.Range("E34").Value = '  -1.66%  '

# r35
$ws.Range("E35").Value = '  +0.41%  '

# r36
$ws.Range("E36").Value = '  -0.92%  '

# r37
$ws.Range("E37").Value = '  -2.48%  '

# r38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.895.54'
$ws.Range("E38").Value = '  -5.95%  '

# r39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0751'
$ws.Range("E39").Value = '  -3.03%  '

# r40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.77'
$ws.Range("E40").Value = '  +3.64%  '

# r41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.37'
$ws.Range("E41").Value = '  -3.94%  '

# r42
$ws.Range("E42").Value = '  -0.29%  '

# r43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.24'
$ws.Range("E43").Value = '  -0.13%  '

# r44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0317'
$ws.Range("E44").Value = '  -2.06%  '

# r45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.774'
$ws.Range("E45").Value = '  -1.03%  '

# r46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.90'
$ws.Range("E46").Value = '  -0.27%  '

# r47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.27'
$ws.Range("E47").Value = '  +1.55%  '

# r48
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.08'
$ws.Range("E48").Value = '  -4.00%  '

# r49
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '316.53'
$ws.Range("E49").Value = '  -0.54%  '

# r50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.53'
$ws.Range("E50").Value = '  -3.06%  '

# r51
$ws.Range("E51").Value = '  -3.15%  '
